$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date in column C was bumped by one day
# (serial 46074 -> 46075, i.e. 2026-02-21 -> 2026-02-22) for every
# data row (rows 2 through 548).
$ws.Range("C2:C548").Value = 46075
